# Update the "Password" column (C) values on the "reg" worksheet.
# Each existing password has a single trailing digit appended (1..9,0)
# matching the corresponding row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

$ws.Range("C2").Value  = "akhilbingi32131211"
$ws.Range("C3").Value  = "akhilbingi64234442322"
$ws.Range("C4").Value  = "akhilbingi3244553433"
$ws.Range("C5").Value  = "akhilbingi3215664544"
$ws.Range("C6").Value  = "akhilbingi32127775655"
$ws.Range("C7").Value  = "nffmf88676"
$ws.Range("C8").Value  = "akhilhdiw599787"
$ws.Range("C9").Value  = "akhilhdingi900898"
$ws.Range("C10").Value = "RoyalEnfield19"
$ws.Range("C11").Value = "RoyalEnfield310"
